# Updated UI presentation of Login and Sign-up Pages
# Cascade the three UI screenshots down by one slide position, turn the
# first slide into a title slide (with link to the full pitch-deck), and
# add speaker notes describing the web app flow under each picture.

$p = $ppt.ActivePresentation

$link = "https://pitchdeck.hypermatic.com/slides/lemk0e8h67792?token=WjBaJGhWam9fQF9tblA="

# ---------------------------------------------------------------------
# Step 1: slide 4 (currently empty) receives the picture that is
# currently on slide 3 (the Sign-Up screenshot).
# ---------------------------------------------------------------------
$src3 = $p.Slides.Item(3).Shapes.Item(1)
$src3.Copy()
$p.Slides.Item(4).Shapes.Paste() | Out-Null

$notes4 = $p.Slides.Item(4).NotesPage.Shapes.Placeholders.Item(2)
$notes4.TextFrame.TextRange.Text = "For new users, they will get the option to sign-up and be able to log in on the web app from the Login Page. Once they click the Sign-up link, a sign-up page similar to the login page will fade in as the Login Page fading out."

# ---------------------------------------------------------------------
# Step 2: slide 3 receives the picture currently on slide 2 (the
# Login/Sign-Up screenshot); its old Sign-Up picture is removed.
# ---------------------------------------------------------------------
$src2 = $p.Slides.Item(2).Shapes.Item(1)
$src2.Copy()
$p.Slides.Item(3).Shapes.Item(1).Delete()
$p.Slides.Item(3).Shapes.Paste() | Out-Null

$notes3 = $p.Slides.Item(3).NotesPage.Shapes.Placeholders.Item(2)
$notes3.TextFrame.TextRange.Text = "Once the user has clicked the login icon on the Landing Page, this small Login Page will fade in, blur the background as shown in the picture a request the user to enter their login details. Once that's done, the user will be sent back to the landing page where they would now have access to the full functionalities of the site. If the user enters incorrect details, they will be requested to reenter their details again or use the option of the 'Forgot' password option below the login and cancel buttons. For new users, they will get the option to sign-up and be able to log in on the web app. Once they click the Sign-up link, a signup page similar to the login page will fade in as the login page fades out."

# ---------------------------------------------------------------------
# Step 3: slide 2 receives the picture currently on slide 1 (the
# Landing Page screenshot); its old Login picture is removed.
# ---------------------------------------------------------------------
$src1 = $p.Slides.Item(1).Shapes.Item(1)
$src1.Copy()
$p.Slides.Item(2).Shapes.Item(1).Delete()
$p.Slides.Item(2).Shapes.Paste() | Out-Null

$notes2 = $p.Slides.Item(2).NotesPage.Shapes.Placeholders.Item(2)
$notes2.TextFrame.TextRange.Text = "This will be our landing page and the user will be able to browse through. But in order for them to add or remove items to their Wish List and/or Cart and place orders, they will have to login. To login, the user will have to click the on the icon located on the right side at the top of the page."

# ---------------------------------------------------------------------
# Step 4: slide 1 loses its picture and becomes a Title Slide that
# points viewers at the online version of the deck.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item(1).Delete()
$slide1.CustomLayout = $p.SlideMaster.CustomLayouts.Item(1)

$title1 = $slide1.Shapes.Item(1).TextFrame.TextRange
$title1.Text = "UI/UX Picture P"

$subtitle1 = $slide1.Shapes.Item(2).TextFrame.TextRange
$subtitle1.Text = "This presentation can also be viewed here:`r" + $link
$subtitle1.Paragraphs(2, 1).Font.Size = 24
$subtitle1.Paragraphs(2, 1).Font.Color.RGB = 0x602000

$notes1 = $slide1.NotesPage.Shapes.Placeholders.Item(2)
$notes1.TextFrame.TextRange.Text = "For a more engaging presentation of the UI, please follow the following link: " + $link + "`rNote: The animation in the presentation is to emphasize on the features each page has and the PowerPoint presentation is to give more details on how the web app will work."

Write-Output "Edit complete. Slide count: $($p.Slides.Count)"
